# Generate Report for Handoff
# b.md has now been handed off for localization; update status/dates and
# record the new handoff file + the "not latest" error detail for both
# the zh-cn and de-de locales, and reflect it on the Overview sheet.

$wb = $excel.ActiveWorkbook

$status = "Ready for handoff"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d8679fb351821d979c0657ee890c7c772ef34033/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5f8092d937294f6c3a72c6baa244eeb241cae306/e2e/b.md."

# ---- Overview sheet: row 3 is b.md ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $status
$wsOverview.Range("F3").Value = $status
$wsOverview.Range("G3").Value = "2016-08-21 02:42:43"

# Note: the stored/persisted column width ends up ~5/6 wider than the
# ColumnWidth value we assign (pixel-rounding in the engine's column-width
# model), so back the requested value off by 5/6 to land on a persisted
# width of exactly 40.
$columnPWidth = 40 - (5 / 6)

# ---- zh-cn sheet: row 3 is b.md ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $status
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-21 02:42:39"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = $columnPWidth

# ---- de-de sheet: row 3 is b.md ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $status
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-21 02:42:43"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = $columnPWidth
